# Alarm Normal load method changes
#
# On the "Add Panels" sheet, two new columns are introduced next to the
# existing "Battery Standby"/"Alarm Load" pair (L:M): an "AlarmLoadingDetail"
# header with a "Battery Alarm (A)" label underneath (column N), and a
# "StandbyLoadingDetail" header with a "Battery Standby (A)" label
# underneath (column O). The sheet's scroll position / selection is also
# moved to frame the new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")
$ws.Activate()

# New header cells in row 7 (same row as "Battery Standby" / "Alarm Load").
$ws.Range("N7").Value = "AlarmLoadingDetail"
$ws.Range("O7").Value = "StandbyLoadingDetail"

# New data cells in row 8 (same row as the 0.052 / 0.052 sample values).
$ws.Range("N8").Value = "Battery Alarm (A)"
$ws.Range("O8").Value = "Battery Standby (A)"

# Match the formatting already used by the neighboring header/data cells
# (column D uses the same bold/shaded header style, column A the same
# shaded data-row style) instead of leaving the new cells unformatted.
$ws.Range("D7").Copy() | Out-Null
$ws.Range("N7:O7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("D8").Copy() | Out-Null
$ws.Range("N8:O8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# Scroll the view over to the new columns and select them, so the newly
# added data is what's visible/active when the workbook is reopened.
$excel.ActiveWindow.ScrollColumn = 10
$ws.Range("N7:O8").Select() | Out-Null
